$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifts existing rows 2-22 down to 3-23)
$ws.Rows.Item(2).Insert()

# Set the new accelerometer reading in the newly inserted row 2
$ws.Cells.Item(2, 1).Value = -3.555192089080811
$ws.Cells.Item(2, 2).Value = 4.907798504829406
$ws.Cells.Item(2, 3).Value = -2.923101136088372

# Remove the two oldest rows, which are now rows 22-23 after the insert,
# keeping the sliding window at 20 data rows (rows 2-21)
$ws.Range("A22:C23").EntireRow.Delete()
